$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the "Price" (column D) text changed. These values must stay as
# TEXT (matching the original inline-string cells), not be auto-converted to
# numbers by Excel's type inference. We force text by temporarily applying a
# text number format, then clear the formatting again so no stray style is
# left behind on the cell.
$priceUpdates = @{
    2  = "40.120.51"
    3  = "2.219.89"
    5  = "293.77"
    6  = "87.56"
    10 = "30.65"
    11 = "50.87"
    14 = "6.41"
    15 = "2.564.55"
    16 = "13.80"
    17 = "2.238.74"
    18 = "0.735"
    19 = "40.060.03"
    20 = "0.0₃0889"
    21 = "11.26"
    22 = "5.77"
    23 = "65.62"
    24 = "235.91"
    28 = "23.12"
    29 = "9.32"
    31 = "158.91"
    32 = "31.77"
    36 = "0.0713"
    39 = "1.76"
    40 = "0.0996"
    41 = "15.61"
    42 = "2.078.66"
    43 = "3.75"
    44 = "19.25"
    46 = "10.01"
    48 = "1.91"
    49 = "2.438.46"
    51 = "1.47"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.ClearFormats()
}

# Rows where the "Volume(1h)" (column E) text changed. These are percentage
# strings padded with spaces, so Excel's type inference already leaves them
# as plain text -- a normal Value assignment is safe.
$volumeUpdates = @{
    2  = "  +0.63%  "
    3  = "  +0.62%  "
    4  = "  +0.04%  "
    5  = "  +1.89%  "
    6  = "  +0.50%  "
    7  = "  -0.06%  "
    8  = "  -0.07%  "
    10 = "  +1.07%  "
    11 = "  +7.18%  "
    12 = "  +0.69%  "
    13 = "  +3.70%  "
    14 = "  -0.14%  "
    15 = "  +0.55%  "
    16 = "  -0.58%  "
    17 = "  +1.41%  "
    18 = "  +1.45%  "
    19 = "  +0.64%  "
    20 = "  +1.10%  "
    21 = "  -2.40%  "
    22 = "  -0.01%  "
    23 = "  +0.44%  "
    24 = "  +0.29%  "
    25 = "  +0.03%  "
    26 = "  +1.44%  "
    27 = "  +0.36%  "
    28 = "  +3.11%  "
    29 = "  +1.60%  "
    30 = "  -5.76%  "
    31 = "  +4.15%  "
    32 = "  +0.56%  "
    33 = "  +0.01%  "
    34 = "  +0.81%  "
    35 = "  +7.83%  "
    36 = "  -0.14%  "
    37 = "  -2.12%  "
    38 = "  +1.69%  "
    39 = "  +4.02%  "
    40 = "  +1.27%  "
    41 = "  -0.27%  "
    42 = "  -0.48%  "
    44 = "  +10.02%  "
    45 = "  +1.22%  "
    46 = "  +0.73%  "
    47 = "  +4.46%  "
    48 = "  -10.80%  "
    49 = "  +0.62%  "
    50 = "  +4.34%  "
    51 = "  +2.12%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
